$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the saved selection/view on the two existing sheets that changed.
# ---------------------------------------------------------------------------
$wsCross = $wb.Worksheets.Item("Cross_valid")
$wsCross.Activate()
$wsCross.Range("V102:W102").Select()

$wsEmp = $wb.Worksheets.Item("Empiric_test")
$wsEmp.Activate()
$wsEmp.Range("V8:V11").Select()

# ---------------------------------------------------------------------------
# 2) Add the new "Sheet1" worksheet at the end of the workbook.
# ---------------------------------------------------------------------------
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet1"

# ---------------------------------------------------------------------------
# 3) Fill in the three stacked SNR / Correct% / Wrong% tables.
# ---------------------------------------------------------------------------

# --- Block 1: rows 2-6 (B = formula 1-C) ------------------------------------
$ws.Range("A2").Value = "SNR"
$ws.Range("B2").Value = "Correct [%]"
$ws.Range("C2").Value = "Wrong [%]"

$ws.Range("A3").Value = 0
$ws.Range("C3").Value = 0.376
$ws.Range("B3").Formula = "=1-C3"

$ws.Range("A4").Value = 6
$ws.Range("C4").Value = 0.347
$ws.Range("B4").Formula = "=1-C4"

$ws.Range("A5").Value = 12
$ws.Range("C5").Value = 0.293
$ws.Range("B5").Formula = "=1-C5"

$ws.Range("A6").Value = 18
$ws.Range("C6").Value = 0.265
$ws.Range("B6").Formula = "=1-C6"

# --- Block 2: rows 8-12 (plain values) --------------------------------------
$ws.Range("A8").Value = "SNR"
$ws.Range("B8").Value = "Correct [%]"
$ws.Range("C8").Value = "Wrong [%]"

$ws.Range("A9").Value = 0
$ws.Range("B9").Value = 0.896
$ws.Range("C9").Value = 0.104

$ws.Range("A10").Value = 6
$ws.Range("B10").Value = 0.963
$ws.Range("C10").Value = 0.037

$ws.Range("A11").Value = 12
$ws.Range("B11").Value = 0.989
$ws.Range("C11").Value = 0.011

$ws.Range("A12").Value = 18
$ws.Range("B12").Value = 0.99
$ws.Range("C12").Value = 0.005

# --- Block 3: rows 14-18 (C = formula 1-B) ----------------------------------
$ws.Range("A14").Value = "SNR"
$ws.Range("B14").Value = "Correct [%]"
$ws.Range("C14").Value = "Wrong [%]"

$ws.Range("A15").Value = 0
$ws.Range("B15").Value = 0.892
$ws.Range("C15").Formula = "=1-B15"

$ws.Range("A16").Value = 6
$ws.Range("B16").Value = 0.911
$ws.Range("C16").Formula = "=1-B16"

$ws.Range("A17").Value = 12
$ws.Range("B17").Value = 0.927
$ws.Range("C17").Formula = "=1-B17"

$ws.Range("A18").Value = 18
$ws.Range("B18").Value = 0.947
$ws.Range("C18").Formula = "=1-B18"

# ---------------------------------------------------------------------------
# 4) Formatting: thin box border around every A:C row used, percent number
#    format on the Correct/Wrong data columns.
# ---------------------------------------------------------------------------
$ws.Range("A2:C6").Borders.LineStyle = 1
$ws.Range("A8:C12").Borders.LineStyle = 1
$ws.Range("A14:C18").Borders.LineStyle = 1

$ws.Range("B3:C6").NumberFormat = "0%"
$ws.Range("B9:C12").NumberFormat = "0%"
$ws.Range("B15:C18").NumberFormat = "0%"

# ---------------------------------------------------------------------------
# 5) Final selection/active state on the new sheet.
# ---------------------------------------------------------------------------
$ws.Range("N19").Select()

Write-Output "done"
